# Apply cell updates to the 'cryptos' worksheet (Coin / Link / Price / Volume(1h))
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.933.93'
$ws.Range("E2").Value = '  +7.52%  '

# Row 3
$ws.Range("D3").Value = '2.647.46'
$ws.Range("E3").Value = '  +9.06%  '

# Row 4
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.31%  '

# Row 5
$ws.Range("D5").Value = '''511.49'
$ws.Range("E5").Value = '  +5.12%  '

# Row 6
$ws.Range("D6").Value = '''157.29'
$ws.Range("E6").Value = '  +2.85%  '

# Row 7
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.04%  '

# Row 8
$ws.Range("D8").Value = '''0.604'
$ws.Range("E8").Value = '  +0.11%  '

# Row 9
$ws.Range("D9").Value = '2.667.37'
$ws.Range("E9").Value = '  +9.91%  '

# Row 10
$ws.Range("D10").Value = '''6.39'
$ws.Range("E10").Value = '  +11.58%  '

# Row 11
$ws.Range("E11").Value = '  +5.48%  '

# Row 12
$ws.Range("D12").Value = '''0.348'
$ws.Range("E12").Value = '  +4.21%  '

# Row 13
$ws.Range("E13").Value = '  +1.22%  '

# Row 14
$ws.Range("D14").Value = '3.130.73'
$ws.Range("E14").Value = '  +10.19%  '

# Row 15
$ws.Range("D15").Value = '60.968.30'
$ws.Range("E15").Value = '  +7.09%  '

# Row 16
$ws.Range("D16").Value = '''21.79'
$ws.Range("E16").Value = '  +5.11%  '

# Row 17
$ws.Range("E17").Value = '  +5.08%  '

# Row 18
$ws.Range("D18").Value = '2.666.47'
$ws.Range("E18").Value = '  +9.60%  '

# Row 19
$ws.Range("D19").Value = '''4.80'
$ws.Range("E19").Value = '  +1.08%  '

# Row 20
$ws.Range("D20").Value = '''348.86'
$ws.Range("E20").Value = '  +7.74%  '

# Row 21
$ws.Range("D21").Value = '''10.51'
$ws.Range("E21").Value = '  +5.59%  '

# Row 22
$ws.Range("D22").Value = '''6.18'
$ws.Range("E22").Value = '  +4.11%  '

# Row 23
$ws.Range("E23").Value = '  -0.01%  '

# Row 24
$ws.Range("D24").Value = '''60.22'
$ws.Range("E24").Value = '  +3.92%  '

# Row 25
$ws.Range("D25").Value = '''0.422'
$ws.Range("E25").Value = '  +3.49%  '

# Row 26
$ws.Range("D26").Value = '2.766.91'
$ws.Range("E26").Value = '  +9.81%  '

# Row 27
$ws.Range("E27").Value = '  +3.80%  '

# Row 28
$ws.Range("D28").Value = '''0.998'
$ws.Range("E28").Value = '  +0.04%  '

# Row 29
$ws.Range("E29").Value = '  +9.70%  '

# Row 30
$ws.Range("D30").Value = '''7.54'
$ws.Range("E30").Value = '  +3.38%  '

# Row 31
$ws.Range("E31").Value = '  +0.04%  '

# Row 32
$ws.Range("D32").Value = '''157.84'
$ws.Range("E32").Value = '  +5.46%  '

# Row 33
$ws.Range("D33").Value = '''19.51'
$ws.Range("E33").Value = '  +5.18%  '

# Row 34
$ws.Range("D34").Value = '''1.58'
$ws.Range("E34").Value = '  +3.69%  '

# Row 35
$ws.Range("D35").Value = '''5.71'
$ws.Range("E35").Value = '  +6.83%  '

# Row 36
$ws.Range("D36").Value = '''4.04'
$ws.Range("E36").Value = '  +8.98%  '

# Row 37
$ws.Range("D37").Value = '''1.22'
$ws.Range("E37").Value = '  +5.99%  '

# Row 38
$ws.Range("E38").Value = '  +11.34%  '

# Row 39
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = '''310.44'
$ws.Range("E39").Value = '  +15.77%  '

# Row 40
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").Value = '''0.860'
$ws.Range("E40").Value = '  +1.88%  '

# Row 41
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '''3.77'
$ws.Range("E41").Value = '  +6.86%  '

# Row 42
$ws.Range("B42").Value = 'SuiNetwork'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D42").Value = '''0.838'
$ws.Range("E42").Value = '  +28.15%  '

# Row 43
$ws.Range("D43").Value = '''35.38'
$ws.Range("E43").Value = '  +3.67%  '

# Row 44
$ws.Range("D44").Value = '''0.643'
$ws.Range("E44").Value = '  +8.73%  '

# Row 45
$ws.Range("E45").Value = '  +8.92%  '

# Row 46
$ws.Range("E46").Value = '  -0.73%  '

# Row 47
$ws.Range("D47").Value = '''0.997'
$ws.Range("E47").Value = '  -0.07%  '

# Row 48
$ws.Range("E48").Value = '  +14.67%  '

# Row 49
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '''0.0237'
$ws.Range("E49").Value = '  +3.93%  '

# Row 50
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '''4.83'
$ws.Range("E50").Value = '  +5.83%  '

# Row 51
$ws.Range("D51").Value = '2.042.19'
$ws.Range("E51").Value = '  +9.29%  '
